# Calculadora.xlsx - BetaPhase 0.1.1 update
# Update the input values that drive the pixel/percentage conversion
# calculator and move the active selection down to the next blank
# separator row (A8:K8), mirroring what was captured in the source
# control diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# "Height" input (H3) changes from 2130 to 713; F3 (=H3/E3) and the
# dependent formulas in G5, E7 and G7 recalculate automatically.
$ws.Range("H3").Value = 713

# "Em Porcentagem" input (C7) changes from 30% to 40%.
$ws.Range("C7").Value = 0.4

# Move the active selection from the A6:K6 spacer row to the A8:K8
# spacer row.
$null = $ws.Range("A8:K8").Select()
